$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.940.80'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '2.405.75'
$ws.Range('E3').Value = '  +1.79%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'" + '553.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').Value = "'" + '142.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.66%  '
$ws.Range('D8').Value = "'" + '0.529'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('D9').Value = '2.398.88'
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').Value = "'" + '25.94'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.21%  '
$ws.Range('E15').Value = '  +5.49%  '
$ws.Range('D16').Value = '2.840.42'
$ws.Range('E16').Value = '  +2.42%  '
$ws.Range('D17').Value = '61.953.91'
$ws.Range('E17').Value = '  +1.75%  '
$ws.Range('D18').Value = '2.403.11'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('E19').Value = '  +3.34%  '
$ws.Range('D20').Value = "'" + '4.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('D21').Value = "'" + '322.66'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('D22').Value = "'" + '6.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').Value = "'" + '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = "'" + '65.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.71%  '
$ws.Range('D25').Value = "'" + '1.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.30%  '
$ws.Range('D26').Value = "'" + '8.99'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.24%  '
$ws.Range('D27').Value = "'" + '573.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +16.02%  '
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').Value = '2.524.61'
$ws.Range('E29').Value = '  +2.03%  '
$ws.Range('D30').Value = "'" + '8.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.14%  '
$ws.Range('D31').Value = '0.0₃0923'
$ws.Range('E31').Value = '  +5.71%  '
$ws.Range('D32').Value = "'" + '1.44'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.79%  '
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('E34').Value = '  +2.55%  '
$ws.Range('D35').Value = "'" + '1.55'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.79%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').Value = "'" + '5.62'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.28%  '
$ws.Range('E39').Value = '  +1.18%  '
$ws.Range('D40').Value = "'" + '150.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.73%  '
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = "'" + '2.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +14.33%  '
$ws.Range('D45').Value = "'" + '149.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.04%  '
$ws.Range('D47').Value = "'" + '0.0538'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.94%  '
$ws.Range('D48').Value = "'" + '20.06'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.04%  '
$ws.Range('D49').Value = "'" + '0.586'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.56%  '
$ws.Range('D50').Value = "'" + '0.0923'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.05%  '
$ws.Range('E51').Value = '  +2.52%  '
